$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.568.52"
$ws.Range("E2").Value = "  -1.99%  "
$ws.Range("D3").Value = "2.293.55"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'302.93"
$ws.Range("E5").Value = "  -2.62%  "
$ws.Range("D6").Value = "'97.97"
$ws.Range("E6").Value = "  -6.02%  "
$ws.Range("D7").Value = "'0.503"
$ws.Range("E7").Value = "  -5.75%  "
$ws.Range("E8").Value = "  +0.14%  "
$ws.Range("D9").Value = "'0.499"
$ws.Range("E9").Value = "  -5.82%  "
$ws.Range("D10").Value = "'34.17"
$ws.Range("E10").Value = "  -6.95%  "
$ws.Range("D11").Value = "'51.46"
$ws.Range("E11").Value = "  -2.73%  "
$ws.Range("D12").Value = "'0.0786"
$ws.Range("E12").Value = "  -3.50%  "
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("D14").Value = "'6.69"
$ws.Range("E14").Value = "  -4.70%  "
$ws.Range("D15").Value = "2.653.16"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("D16").Value = "'15.53"
$ws.Range("E16").Value = "  +2.56%  "
$ws.Range("D17").Value = "2.323.72"
$ws.Range("E17").Value = "  +0.49%  "
$ws.Range("D18").Value = "'0.804"
$ws.Range("E18").Value = "  -0.91%  "
$ws.Range("D19").Value = "42.553.19"
$ws.Range("E19").Value = "  -1.82%  "
$ws.Range("D20").Value = "0.0₃0896"
$ws.Range("E20").Value = "  -3.43%  "
$ws.Range("D21").Value = "'11.43"
$ws.Range("E21").Value = "  -6.12%  "
$ws.Range("D22").Value = "'6.03"
$ws.Range("E22").Value = "  -2.49%  "
$ws.Range("D23").Value = "'68.50"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D24").Value = "'233.89"
$ws.Range("E24").Value = "  -3.67%  "
$ws.Range("D25").Value = "'1.96"
$ws.Range("E25").Value = "  -3.36%  "
$ws.Range("D26").Value = "'2.51"
$ws.Range("E26").Value = "  -4.11%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("D28").Value = "'24.90"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("E29").Value = "  -7.58%  "
$ws.Range("D30").Value = "'34.45"
$ws.Range("E30").Value = "  -6.95%  "
$ws.Range("D31").Value = "'9.13"
$ws.Range("E31").Value = "  -5.25%  "
$ws.Range("D32").Value = "'162.74"
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("D34").Value = "'5.00"
$ws.Range("E34").Value = "  -5.67%  "
$ws.Range("D35").Value = "'4.60"
$ws.Range("E35").Value = "  +2.75%  "
$ws.Range("E36").Value = "  -3.75%  "
$ws.Range("D37").Value = "'0.0710"
$ws.Range("E37").Value = "  -4.53%  "
$ws.Range("D38").Value = "'16.87"
$ws.Range("E38").Value = "  -8.29%  "
$ws.Range("D39").Value = "'2.87"
$ws.Range("E39").Value = "  -6.45%  "
$ws.Range("E40").Value = "  -4.95%  "
$ws.Range("E41").Value = "  -5.80%  "
$ws.Range("E42").Value = "  -4.35%  "
$ws.Range("E43").Value = "  -9.53%  "
$ws.Range("D44").Value = "1.979.56"
$ws.Range("E44").Value = "  -0.66%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0278"
$ws.Range("E45").Value = "  -5.10%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'18.54"
$ws.Range("E46").Value = "  -3.34%  "
$ws.Range("D47").Value = "'10.11"
$ws.Range("E47").Value = "  +0.97%  "
$ws.Range("D48").Value = "'2.85"
$ws.Range("E48").Value = "  -7.09%  "
$ws.Range("D49").Value = "'55.02"
$ws.Range("E49").Value = "  -2.15%  "
$ws.Range("E50").Value = "  -2.92%  "
$ws.Range("D51").Value = "2.524.45"
$ws.Range("E51").Value = "  -0.54%  "
